# Generate Report for Handoff
#
# The localization-status report is regenerated: rows 7 and 10-16 (the
# "Ready for handoff" / "Handback transform failed" rows whose handback
# had not actually completed, timestamp 0001-01-01) get their
# "Latest Handoff Datetime" stamped with a single, fresh run timestamp on
# each sheet (Overview uses column D, the per-language sheets use column E).

$wb = $excel.ActiveWorkbook

# --- Overview sheet: column D ("Latest Handoff Date") ---
$overview = $wb.Worksheets.Item("Overview")
$overviewRows = 7,10,11,12,13,14,15,16
foreach ($r in $overviewRows) {
    $overview.Cells.Item($r, 4).Value = "2016-28-21 04:28:09"
}

# --- zh-cn sheet: column E ("Latest Handoff Datetime") ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcnRows = 7,10,11,12,13,14,15,16
foreach ($r in $zhcnRows) {
    $zhcn.Cells.Item($r, 5).Value = "2016-03-21 04:28:04"
}

# --- de-de sheet: column E ("Latest Handoff Datetime") ---
$dede = $wb.Worksheets.Item("de-de")
$dedeRows = 7,10,11,12,13,14,15,16
foreach ($r in $dedeRows) {
    $dede.Cells.Item($r, 5).Value = "2016-03-21 04:28:09"
}
